$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17: H17,J17,L17,N17
$ws.Range("H17").Value = 1956.6666
$ws.Range("J17").Value = 1956.6666
$ws.Range("L17").Value = 5869.9998
$ws.Range("N17").Value = -6205.9998
# Row 40: H40,I40,J40,K40,L40,M40,N40
$ws.Range("H40").Value = 2381.25
$ws.Range("I40").Value = 2100
$ws.Range("J40").Value = 2550
$ws.Range("K40").Value = 2100
$ws.Range("L40").Value = 2550
$ws.Range("M40").Value = -1925
$ws.Range("N40").Value = -2900
# Row 138: H138,J138,L138,N138
$ws.Range("H138").Value = 4559.923
$ws.Range("J138").Value = 5970.5483
$ws.Range("L138").Value = 17911.6449
$ws.Range("N138").Value = -28191.6449

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45: H45,I45,K45,M45
$ws.Range("H45").Value = 2054.5
$ws.Range("I45").Value = 1998.6666
$ws.Range("K45").Value = 1998.6666
$ws.Range("M45").Value = -1621.6666
# Row 74: H74,I74,J74,K74,L74,M74,N74
$ws.Range("H74").Value = 1140
$ws.Range("I74").Value = 1180
$ws.Range("J74").Value = 740
$ws.Range("K74").Value = 1180
$ws.Range("L74").Value = 740
$ws.Range("M74").Value = -306
$ws.Range("N74").Value = -2488
# Row 77: H77,I77,J77,K77,L77,M77,N77
$ws.Range("H77").Value = 1140
$ws.Range("I77").Value = 1180
$ws.Range("J77").Value = 740
$ws.Range("K77").Value = 5900
$ws.Range("L77").Value = 3700
$ws.Range("M77").Value = -1532
$ws.Range("N77").Value = -12436
# Row 80: H80,I80,J80,K80,L80,M80,N80
$ws.Range("H80").Value = 45062.2
$ws.Range("I80").Value = 25100
$ws.Range("J80").Value = 50052.75
$ws.Range("K80").Value = 25100
$ws.Range("L80").Value = 50052.75
$ws.Range("M80").Value = -24102
$ws.Range("N80").Value = -52048.75
# Row 83: H83,I83,J83,K83,L83,M83,N83
$ws.Range("H83").Value = 45062.2
$ws.Range("I83").Value = 25100
$ws.Range("J83").Value = 50052.75
$ws.Range("K83").Value = 75300
$ws.Range("L83").Value = 150158.25
$ws.Range("M83").Value = -70308
$ws.Range("N83").Value = -160142.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 32: H32,J32,L32,N32
$ws.Range("H32").Value = 1000
$ws.Range("J32").Value = 1000
$ws.Range("L32").Value = 1000
$ws.Range("N32").Value = -1768
# Row 48: H48,J48,L48,N48
$ws.Range("H48").Value = 77777
$ws.Range("J48").Value = 77777
$ws.Range("L48").Value = 77777
$ws.Range("N48").Value = -78607
# Row 86: H86,I86,J86,K86,L86,M86,N86
$ws.Range("H86").Value = 3963.3333
$ws.Range("I86").Value = 3963.3333
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3963.3333
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -2840.3333
# Row 89: H89,I89,J89,K89,L89,M89,N89
$ws.Range("H89").Value = 3963.3333
$ws.Range("I89").Value = 3963.3333
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 19816.6665
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -14200.6665
# Row 105: H105,I105,K105,M105
$ws.Range("H105").Value = 3066.7856
$ws.Range("I105").Value = 3066.7856
$ws.Range("K105").Value = 3066.7856
$ws.Range("M105").Value = -1319.7856

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 12: H12,I12,K12,M12
$ws.Range("H12").Value = 4790.125
$ws.Range("I12").Value = 2053.6667
$ws.Range("K12").Value = 2053.6667
$ws.Range("M12").Value = -1883.6667
# Row 16: H16,I16,J16,K16,L16,M16,N16
$ws.Range("H16").Value = 4722
$ws.Range("I16").Value = 4722
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 4722
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -4435
# Row 35: H35,I35,J35,K35,L35,M35,N35
$ws.Range("H35").Value = 1356.25
$ws.Range("I35").Value = 1516.6666
$ws.Range("J35").Value = 875
$ws.Range("K35").Value = 1516.6666
$ws.Range("L35").Value = 875
$ws.Range("M35").Value = -1222.6666
$ws.Range("N35").Value = -1463
# Row 58: H58,I58,K58,M58
$ws.Range("H58").Value = 2311.3333
$ws.Range("I58").Value = 1298
$ws.Range("K58").Value = 1298
$ws.Range("M58").Value = -1095
# Row 86: H86,I86,K86,M86
$ws.Range("H86").Value = 7248.75
$ws.Range("I86").Value = 4248.5
$ws.Range("K86").Value = 4248.5
$ws.Range("M86").Value = -3125.5
# Row 89: H89,I89,K89,M89
$ws.Range("H89").Value = 7248.75
$ws.Range("I89").Value = 4248.5
$ws.Range("K89").Value = 21242.5
$ws.Range("M89").Value = -15626.5
# Row 99: H99,I99,K99,M99
$ws.Range("H99").Value = 12398.087
$ws.Range("I99").Value = 9167.308000000001
$ws.Range("K99").Value = 9167.308000000001
$ws.Range("M99").Value = -7669.308000000001
# Row 113: H113,I113,J113,K113,L113,M113,N113
$ws.Range("H113").Value = 4722
$ws.Range("I113").Value = 4722
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4722
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -2552
# Row 126: H126,I126,K126,M126
$ws.Range("H126").Value = 12398.087
$ws.Range("I126").Value = 9167.308000000001
$ws.Range("K126").Value = 27501.924
$ws.Range("M126").Value = -25031.924
# Row 136: H136,I136,K136,M136
$ws.Range("H136").Value = 2311.3333
$ws.Range("I136").Value = 1298
$ws.Range("K136").Value = 3894
$ws.Range("M136").Value = -1344

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 23: H23,I23,J23,K23,L23,M23,N23
$ws.Range("H23").Value = 120.625
$ws.Range("I23").Value = 155.875
$ws.Range("J23").Value = 85.375
$ws.Range("K23").Value = 467.625
$ws.Range("L23").Value = 256.125
$ws.Range("M23").Value = -232.625
$ws.Range("N23").Value = -726.125
# Row 34: H34,J34,L34,N34
$ws.Range("H34").Value = 882
$ws.Range("J34").Value = 4000
$ws.Range("L34").Value = 12000
$ws.Range("N34").Value = -12168
# Row 39: H39,J39,L39,N39
$ws.Range("H39").Value = 7666.6665
$ws.Range("J39").Value = 10500
$ws.Range("L39").Value = 31500
$ws.Range("N39").Value = -32088
# Row 55: H55,J55,L55,N55
$ws.Range("H55").Value = 145410.72
$ws.Range("J55").Value = 2866.6667
$ws.Range("L55").Value = 8600.000100000001
$ws.Range("N55").Value = -8954.000100000001
# Row 113: H113,I113,K113,M113
$ws.Range("H113").Value = 3407.5715
$ws.Range("I113").Value = 4899
$ws.Range("K113").Value = 14697
$ws.Range("M113").Value = -12527
# Row 137: H137,I137,J137,K137,L137,M137,N137
$ws.Range("H137").Value = 3275.4
$ws.Range("I137").Value = 3275.4
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 9826.200000000001
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -4726.200000000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 11: H11,I11,J11,K11,L11,M11,N11
$ws.Range("H11").Value = 2207275
$ws.Range("I11").Value = 2867250
$ws.Range("J11").Value = 667333.3
$ws.Range("K11").Value = 2867250
$ws.Range("L11").Value = 667333.3
$ws.Range("M11").Value = -2867111
$ws.Range("N11").Value = -667611.3
# Row 12: H12,I12,K12,M12
$ws.Range("H12").Value = 4300
$ws.Range("I12").Value = 4300
$ws.Range("K12").Value = 4300
$ws.Range("M12").Value = -4160
# Row 41: H41,J41,L41,N41
$ws.Range("H41").Value = 662.75
$ws.Range("J41").Value = 700
$ws.Range("L41").Value = 700
$ws.Range("N41").Value = -1410
# Row 123: H123,J123,L123,N123
$ws.Range("H123").Value = 28210.3
$ws.Range("J123").Value = 28210.3
$ws.Range("L123").Value = 28210.3
$ws.Range("N123").Value = -33110.3
# Row 132: H132,J132,L132,N132
$ws.Range("H132").Value = 1982.2413
$ws.Range("J132").Value = 2366.3333
$ws.Range("L132").Value = 7098.999899999999
$ws.Range("N132").Value = -12158.9999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 136: H136,I136,K136,M136
$ws.Range("H136").Value = 2372.2632
$ws.Range("I136").Value = 2372.2632
$ws.Range("K136").Value = 7116.7896
$ws.Range("M136").Value = -4566.7896

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 75: H75,I75,K75,M75
$ws.Range("H75").Value = 40118
$ws.Range("I75").Value = 40118
$ws.Range("K75").Value = 40118
$ws.Range("M75").Value = -39182
# Row 78: H78,I78,K78,M78
$ws.Range("H78").Value = 40118
$ws.Range("I78").Value = 40118
$ws.Range("K78").Value = 120354
$ws.Range("M78").Value = -115674
# Row 132: H132,J132,L132,N132
$ws.Range("H132").Value = 885.68085
$ws.Range("J132").Value = 1641
$ws.Range("L132").Value = 4923
$ws.Range("N132").Value = -9983
